# Generate Report for Handoff
# Adds two new localization entries (b354880c... and cd9fc48a...) as new
# rows (6 and 7) on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---- shared literal values -------------------------------------------------
$guid1      = "b354880c-24dd-4b16-ac5e-af6f3853cae3"
$hash1      = "c27ecce420d5d769275bc64fc40f150097be80d3"
$guid2      = "cd9fc48a-08cd-4161-afee-beef6704f90c"
$hash2      = "3e289fc56ec5a62f889a45e3fe7009d72ac149fe"

$md1        = "$guid1.md"
$zhcn1      = "$guid1.$hash1.zh-cn.xlf"
$dede1      = "$guid1.$hash1.de-de.xlf"

$md2        = "$guid2.md"
$zhcn2      = "$guid2.$hash2.zh-cn.xlf"
$dede2      = "$guid2.$hash2.de-de.xlf"

$handoffDate1 = "2016-03-24 22:42:03"
$handoffDate2 = "2016-03-24 22:42:03"
$zhHandoffDt  = "2016-03-24 22:41:58"
$deHandoffDt  = "2016-03-24 22:42:03"
$status       = "Ready for handoff"
$ext          = ".md"
$epoch        = "0001-01-01 00:00:00"
$reason       = "Include"

$md1Url   = "https://github.com/OpenLocalizationTest/oltest/blob/688d0611bd4c0c191af6c50c8d39911c55748fb7/e2e/$md1"
$md2Url   = "https://github.com/OpenLocalizationTest/oltest/blob/ebec6f6002e3504078db1150b241869c6c5c4f1f/e2e/$md2"

$zhcn1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/000c84a71e6e5398fbd689e047c4a227c72df679/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhcn1"
$zhcn2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f2f087679a53c277804e906e957f5aa781306dfc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhcn2"

$dede1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75ce76b0a73944c61752bce24936b9a39251f0b8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$dede1"
$dede2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b1a0ab6fb6fbb9d2f6e10f5e4ff22c2a8c7e4a1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$dede2"

# =============================================================================
# Sheet 1: "Overview"
# =============================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Add($ws.Range("A6"), $md1Url, "", "", $md1)
$ws.Range("B6").Value = $status
$ws.Range("C6").Value = $status
$ws.Range("D6").Value = $handoffDate1
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("A7"), $md2Url, "", "", $md2)
$ws.Range("B7").Value = $status
$ws.Range("C7").Value = $status
$ws.Range("D7").Value = $handoffDate2
$ws.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# =============================================================================
# Sheet 2: "zh-cn"
# =============================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Add($ws.Range("A6"), $md1Url, "", "", $md1)
$ws.Range("B6").Value = $ext
$ws.Range("C6").Value = $status
$ws.Hyperlinks.Add($ws.Range("D6"), $zhcn1Url, "", "", $zhcn1)
$ws.Range("E6").Value = $zhHandoffDt
$ws.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H6").Value = $epoch
$ws.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J6").Value = $reason

$ws.Hyperlinks.Add($ws.Range("A7"), $md2Url, "", "", $md2)
$ws.Range("B7").Value = $ext
$ws.Range("C7").Value = $status
$ws.Hyperlinks.Add($ws.Range("D7"), $zhcn2Url, "", "", $zhcn2)
$ws.Range("E7").Value = $zhHandoffDt
$ws.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H7").Value = $epoch
$ws.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J7").Value = $reason

# =============================================================================
# Sheet 3: "de-de"
# =============================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Add($ws.Range("A6"), $md1Url, "", "", $md1)
$ws.Range("B6").Value = $ext
$ws.Range("C6").Value = $status
$ws.Hyperlinks.Add($ws.Range("D6"), $dede1Url, "", "", $dede1)
$ws.Range("E6").Value = $deHandoffDt
$ws.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H6").Value = $epoch
$ws.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J6").Value = $reason

$ws.Hyperlinks.Add($ws.Range("A7"), $md2Url, "", "", $md2)
$ws.Range("B7").Value = $ext
$ws.Range("C7").Value = $status
$ws.Hyperlinks.Add($ws.Range("D7"), $dede2Url, "", "", $dede2)
$ws.Range("E7").Value = $deHandoffDt
$ws.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H7").Value = $epoch
$ws.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J7").Value = $reason
